$d = $word.ActiveDocument

# Each pair below is a unique, whole run/paragraph text found in the
# 'Subjekty' list and later cross-reference paragraphs. PERSON_N tokens
# with N >= 47 are each decremented by 1 (PERSON_47 was removed/merged).
$replacements = @(
    ,@("[[PERSON_46]] – „k [[PERSON_47]]“", "[[PERSON_46]] – „k [[PERSON_46]]“")
    ,@("[[PERSON_48]] – „se [[PERSON_48]]“", "[[PERSON_47]] – „se [[PERSON_47]]“")
    ,@("[[PERSON_49]] – „u [[PERSON_49]]“", "[[PERSON_48]] – „u [[PERSON_48]]“")
    ,@("[[PERSON_50]] – „o [[PERSON_51]]“", "[[PERSON_49]] – „o [[PERSON_50]]“")
    ,@("[[PERSON_52]] – „s [[PERSON_52]]“", "[[PERSON_51]] – „s [[PERSON_51]]“")
    ,@("[[PERSON_53]] – „k [[PERSON_54]]“", "[[PERSON_52]] – „k [[PERSON_53]]“")
    ,@("[[PERSON_55]] – „od [[PERSON_56]]“", "[[PERSON_54]] – „od [[PERSON_55]]“")
    ,@("[[PERSON_57]] – „s [[PERSON_57]]“", "[[PERSON_56]] – „s [[PERSON_56]]“")
    ,@("[[PERSON_58]] – „u [[PERSON_59]]“", "[[PERSON_57]] – „u [[PERSON_58]]“")
    ,@("[[PERSON_60]] – „o [[PERSON_61]]“", "[[PERSON_59]] – „o [[PERSON_60]]“")
    ,@("[[PERSON_62]] – „k [[PERSON_63]]“", "[[PERSON_61]] – „k [[PERSON_62]]“")
    ,@("V těchto řízeních bylo jednáno např. s [[PERSON_3]], [[PERSON_8]], [[PERSON_36]] či [[PERSON_64]].", "V těchto řízeních bylo jednáno např. s [[PERSON_3]], [[PERSON_8]], [[PERSON_36]] či [[PERSON_63]].")
    ,@("Neurologické testy č. NEU/2025/44119 provedené MUDr. [[PERSON_50]],", "Neurologické testy č. NEU/2025/44119 provedené MUDr. [[PERSON_49]],")
    ,@("Zvláštní pozornost byla věnována výsledkům [[PERSON_22]], [[PERSON_28]] a [[PERSON_60]].", "Zvláštní pozornost byla věnována výsledkům [[PERSON_22]], [[PERSON_28]] a [[PERSON_59]].")
    ,@("mobil [[PERSON_65]] S22, [[IMEI_1]],", "mobil [[PERSON_64]] S22, [[IMEI_1]],")
    ,@("[[PERSON_66]] poskytly technické přístupy pro řešení kauz:", "[[PERSON_65]] poskytly technické přístupy pro řešení kauz:")
    ,@("právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_62]]),", "právní cloud účet ID: LEX-ACC-88221 (spravovala [[PERSON_61]]),")
    ,@("[[PERSON_57]] („výslech [[PERSON_57]]“),", "[[PERSON_56]] („výslech [[PERSON_56]]“),")
    ,@("[[PERSON_52]] („výpověď [[PERSON_52]]“),", "[[PERSON_51]] („výpověď [[PERSON_51]]“),")
    ,@("[[PERSON_49]] („záznam o výslechu [[PERSON_49]]“),", "[[PERSON_48]] („záznam o výslechu [[PERSON_48]]“),")
    ,@("[[PERSON_27]] („výslech [[PERSON_67]]“).", "[[PERSON_27]] („výslech [[PERSON_66]]“).")
    ,@("Tyto účty byly doloženy např. od [[PERSON_30]], [[PERSON_53]] nebo [[PERSON_68]].", "Tyto účty byly doloženy např. od [[PERSON_30]], [[PERSON_52]] nebo [[PERSON_67]].")
    ,@("[[PERSON_55]],", "[[PERSON_54]],")
    ,@("[[PERSON_64]],", "[[PERSON_63]],")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $oldText"
    }
}

Write-Output "Done"
